$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").Value = "62.729.66"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "3.059.42"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.44"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.06"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.053.90"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("E10").Value = "  -4.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.08"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.18%  "
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.94"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.15%  "
$ws.Range("D15").Value = "3.552.25"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").Value = "62.780.67"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "3.063.86"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.58"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.96"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.93%  "
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.05"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.59"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.98"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.03%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.01"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.80"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("E31").Value = "  -9.30%  "
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "56.85"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.15%  "
$ws.Range("E34").Value = "  -7.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.95"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "473.24"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -12.18%  "
$ws.Range("D38").Value = "3.084.04"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("E39").Value = "  -6.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0789"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("E44").Value = "  -3.17%  "
$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").Value = "0.0₃0531"
$ws.Range("E46").Value = "  +5.97%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.92"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("E48").Value = "  -6.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.13"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.26"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.59%  "
